$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- column width updates (row of header widths changed) ----
# Target XML widths: 24, 16.7109375, 22.7109375, 19.7109375, 23.7109375, 21, 22.28515625
# This engine quantizes ColumnWidth to a 1/6-character pixel grid (MDW=6), so we pick the
# ColumnWidth input that rounds to the closest achievable grid point to the target width.
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws.Columns.Item(2).ColumnWidth = 15.833333333333334
$ws.Columns.Item(3).ColumnWidth = 21.833333333333332
$ws.Columns.Item(4).ColumnWidth = 18.833333333333332
$ws.Columns.Item(5).ColumnWidth = 22.833333333333332
$ws.Columns.Item(6).ColumnWidth = 20.166666666666668
$ws.Columns.Item(7).ColumnWidth = 21.5

# ---- cell value updates ----
$ws.Range("B2").Value = 0.35013
$ws.Range("C2").Value = 0.16335
$ws.Range("D2").Value = 0.83862999999999999
$ws.Range("E2").Value = 0.77719000000000005
$ws.Range("F2").Value = 0.83496999999999999
$ws.Range("G2").Value = 0.78820000000000001
$ws.Range("B3").Value = 0.04433
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 0.88809000000000005
$ws.Range("E3").Value = 0.89866999999999997
$ws.Range("F3").Value = 0.90403
$ws.Range("G3").Value = 0.91142000000000001
$ws.Range("B4").Value = 0.22517000000000001
$ws.Range("C4").Value = 0.0
$ws.Range("D4").Value = 0.01007
$ws.Range("E4").Value = 0.0
$ws.Range("F4").Value = 0.04271
$ws.Range("G4").Value = 0.0
$ws.Range("B5").Value = 0.34466000000000002
$ws.Range("C5").Value = 0.26367000000000002
$ws.Range("D5").Value = 0.90575000000000006
$ws.Range("E5").Value = 0.91976000000000002
$ws.Range("F5").Value = 0.95084000000000002
$ws.Range("G5").Value = 0.96008000000000004
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 0.0
$ws.Range("D6").Value = 0.00123
$ws.Range("E6").Value = 0.0
$ws.Range("F6").Value = 0.00779
$ws.Range("G6").Value = 0.0
$ws.Range("B7").Value = 0.02532
$ws.Range("C7").Value = 0.0
$ws.Range("D7").Value = 0.00063
$ws.Range("E7").Value = 0.0
$ws.Range("F7").Value = 0.00368
$ws.Range("G7").Value = 0.0
$ws.Range("B8").Value = 0.35224
$ws.Range("C8").Value = 0.0
$ws.Range("D8").Value = 0.84131
$ws.Range("E8").Value = 0.74985999999999997
$ws.Range("F8").Value = 0.88927999999999996
$ws.Range("G8").Value = 0.73614000000000002
$ws.Range("B9").Value = 0.18722
$ws.Range("C9").Value = 0.1149
$ws.Range("D9").Value = 0.75968999999999998
$ws.Range("E9").Value = 0.60409999999999997
$ws.Range("F9").Value = 0.84655999999999998
$ws.Range("G9").Value = 0.69228000000000001
$ws.Range("B10").Value = 0.48475000000000001
$ws.Range("C10").Value = 0.25751000000000002
$ws.Range("D10").Value = 0.89912000000000003
$ws.Range("E10").Value = 0.87217999999999996
$ws.Range("F10").Value = 0.91154000000000002
$ws.Range("G10").Value = 0.87453999999999998
$ws.Range("B11").Value = 0.13428000000000001
$ws.Range("C11").Value = 0.00003
$ws.Range("D11").Value = 0.11205
$ws.Range("E11").Value = 0.00131
$ws.Range("F11").Value = 0.12919
$ws.Range("G11").Value = 0.00148
$ws.Range("B12").Value = 0.0
$ws.Range("C12").Value = 0.0
$ws.Range("D12").Value = 0.00063
$ws.Range("E12").Value = 0.0
$ws.Range("F12").Value = 0.00368
$ws.Range("G12").Value = 0.0
$ws.Range("B13").Value = 0.48475000000000001
$ws.Range("C13").Value = 0.26367000000000002
$ws.Range("D13").Value = 0.90575000000000006
$ws.Range("E13").Value = 0.91976000000000002
$ws.Range("F13").Value = 0.95084000000000002
$ws.Range("G13").Value = 0.96008000000000004
$ws.Range("B14").Value = 0.21940416666666671
$ws.Range("C14").Value = 0.08859416666666665
$ws.Range("D14").Value = 0.5135791666666667
$ws.Range("E14").Value = 0.47856916666666671
$ws.Range("F14").Value = 0.53959250000000014
$ws.Range("G14").Value = 0.49368499999999998
$ws.Range("B15").Value = 0.21940416666666671
$ws.Range("C15").Value = 0.00003
$ws.Range("D15").Value = 0.75968999999999998
$ws.Range("E15").Value = 0.60409999999999997
$ws.Range("F15").Value = 0.83496999999999999
$ws.Range("G15").Value = 0.69228000000000001

# ---- selection change ----
$ws.Range("F15").Select()
